$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = -0.0411
$ws.Cells.Item(2, 5).Value = -0.1515
$ws.Cells.Item(2, 6).Value = 0.1775
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 12691.7
$ws.Cells.Item(2, 12).Value = 0.2646221143581532
$ws.Cells.Item(2, 13).Value = 10974.5
$ws.Cells.Item(2, 14).Value = 0.04011631519883611
$ws.Cells.Item(2, 15).Value = 0.8646989764964503
$ws.Cells.Item(2, 16).Value = 10788.3
$ws.Cells.Item(2, 17).Value = 0.03943567754882716
$ws.Cells.Item(2, 18).Value = 0.8500279710361889
$ws.Cells.Item(2, 19).Value = 186.1999999999996
$ws.Cells.Item(2, 20).Value = 0.01696660440111163
$ws.Cells.Item(2, 21).Value = 117133.4
$ws.Cells.Item(2, 22).Value = 0.4281707954541301
$ws.Cells.Item(2, 23).Value = 0.05573564794305036
$ws.Cells.Item(2, 24).Value = 0.08399779715902757
$ws.Cells.Item(2, 25).Value = -0.02826214921597721
$ws.Cells.Item(2, 26).Value = 0.05969463404713502
$ws.Cells.Item(2, 27).Value = 0
$ws.Cells.Item(2, 28).Value = 0.02968172793475884
$ws.Cells.Item(2, 29).Value = -0.02968172793475884
$ws.Cells.Item(2, 30).Value = 723109.3999999999
$ws.Cells.Item(2, 31).Value = 0
$ws.Cells.Item(2, 32).Value = 723109.3999999999
$ws.Cells.Item(2, 33).Value = 605975.9999999999
$ws.Cells.Item(2, 34).Value = 0.7255207407339032
$ws.Cells.Item(2, 35).Value = 0.7950581035774291
$ws.Cells.Item(2, 36).Value = 0.6889668839385908
$ws.Cells.Item(2, 37).Value = 0.7647622952712722
$ws.Cells.Item(2, 38).Value = 0
$ws.Cells.Item(2, 39).Value = 0
$ws.Cells.Item(2, 40).ClearContents()
$ws.Cells.Item(2, 42).ClearContents()

# Row 3
$ws.Cells.Item(3, 4).Value = -0.0564
$ws.Cells.Item(3, 5).Value = -0.137
$ws.Cells.Item(3, 6).Value = 0.176
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 2564.2
$ws.Cells.Item(3, 12).Value = 0.2400868889450671
$ws.Cells.Item(3, 13).Value = 2138.4
$ws.Cells.Item(3, 14).Value = 0.04322209129110435
$ws.Cells.Item(3, 15).Value = 0.8339443101162157
$ws.Cells.Item(3, 16).Value = 2050.9
$ws.Cells.Item(3, 17).Value = 0.04145351058217635
$ws.Cells.Item(3, 18).Value = 0.799820606816941
$ws.Cells.Item(3, 19).Value = 87.5
$ws.Cells.Item(3, 20).Value = 0.04091844369622148
$ws.Cells.Item(3, 21).Value = 51842.3
$ws.Cells.Item(3, 22).Value = 1.047854762130948
$ws.Cells.Item(3, 23).Value = 0.0625493965088255
$ws.Cells.Item(3, 24).Value = 0.09273701385025866
$ws.Cells.Item(3, 25).Value = -0.03018761734143316
$ws.Cells.Item(3, 26).Value = 0.06011431530423573
$ws.Cells.Item(3, 27).Value = 0
$ws.Cells.Item(3, 28).Value = 0.02960084849567224
$ws.Cells.Item(3, 29).Value = -0.02960084849567224
$ws.Cells.Item(3, 30).Value = 189992.8
$ws.Cells.Item(3, 31).Value = 0
$ws.Cells.Item(3, 32).Value = 189992.8
$ws.Cells.Item(3, 33).Value = 138150.5
$ws.Cells.Item(3, 34).Value = 0.7933970162965747
$ws.Cells.Item(3, 35).Value = 0.8121665449083136
$ws.Cells.Item(3, 36).Value = 0.7363110072634166
$ws.Cells.Item(3, 37).Value = 0.7586893366503561
$ws.Cells.Item(3, 38).Value = 0
$ws.Cells.Item(3, 39).Value = 0
$ws.Cells.Item(3, 40).ClearContents()
$ws.Cells.Item(3, 42).ClearContents()

# Row 4
$ws.Cells.Item(4, 4).Value = -0.04190000000000001
$ws.Cells.Item(4, 5).Value = -0.166
$ws.Cells.Item(4, 6).Value = 0.179
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 1834.4
$ws.Cells.Item(4, 12).Value = 0.1790130083046266
$ws.Cells.Item(4, 13).Value = 1665.2
$ws.Cells.Item(4, 14).Value = 0.02901540681161592
$ws.Cells.Item(4, 15).Value = 0.907762756214566
$ws.Cells.Item(4, 16).Value = 1665.2
$ws.Cells.Item(4, 17).Value = 0.02901540681161592
$ws.Cells.Item(4, 18).Value = 0.907762756214566
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 24840.1
$ws.Cells.Item(4, 22).Value = 0.4328282529072908
$ws.Cells.Item(4, 23).Value = 0.04892189937727522
$ws.Cells.Item(4, 24).Value = 0.08853585395341607
$ws.Cells.Item(4, 25).Value = -0.03961395457614085
$ws.Cells.Item(4, 26).Value = 0.04513811497614758
$ws.Cells.Item(4, 27).Value = 0
$ws.Cells.Item(4, 28).Value = 0.0296359143263922
$ws.Cells.Item(4, 29).Value = -0.0296359143263922
$ws.Cells.Item(4, 30).Value = 205164.9
$ws.Cells.Item(4, 31).Value = 0
$ws.Cells.Item(4, 32).Value = 205164.9
$ws.Cells.Item(4, 33).Value = 180324.8
$ws.Cells.Item(4, 34).Value = 0.7814165483740366
$ws.Cells.Item(4, 35).Value = 0.8236163828143033
$ws.Cells.Item(4, 36).Value = 0.7585756052415707
$ws.Cells.Item(4, 37).Value = 0.8040795068633886
$ws.Cells.Item(4, 38).Value = 0
$ws.Cells.Item(4, 39).Value = 0
$ws.Cells.Item(4, 40).ClearContents()
$ws.Cells.Item(4, 42).ClearContents()

# Row 5
$ws.Cells.Item(5, 4).Value = -0.0403
$ws.Cells.Item(5, 5).Value = -0.222
$ws.Cells.Item(5, 6).Value = 0.189
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 1641.6
$ws.Cells.Item(5, 12).Value = 0.1346677604593929
$ws.Cells.Item(5, 13).Value = 1858.8
$ws.Cells.Item(5, 14).Value = 0.03413667525531801
$ws.Cells.Item(5, 15).Value = 1.132309941520468
$ws.Cells.Item(5, 16).Value = 1805
$ws.Cells.Item(5, 17).Value = 0.03314864366034486
$ws.Cells.Item(5, 18).Value = 1.099537037037037
$ws.Cells.Item(5, 19).Value = 53.79999999999995
$ws.Cells.Item(5, 20).Value = 0.02894340434689044
$ws.Cells.Item(5, 21).Value = 21597.8
$ws.Cells.Item(5, 22).Value = 0.3966414271730727
$ws.Cells.Item(5, 23).Value = 0.03718646010542506
$ws.Cells.Item(5, 24).Value = 0.07945974036463907
$ws.Cells.Item(5, 25).Value = -0.04227328025921401
$ws.Cells.Item(5, 26).Value = 0.06466984162520246
$ws.Cells.Item(5, 27).Value = 0
$ws.Cells.Item(5, 28).Value = 0.02972754154312549
$ws.Cells.Item(5, 29).Value = -0.02972754154312549
$ws.Cells.Item(5, 30).Value = 163452.3
$ws.Cells.Item(5, 31).Value = 0
$ws.Cells.Item(5, 32).Value = 163452.3
$ws.Cells.Item(5, 33).Value = 141854.5
$ws.Cells.Item(5, 34).Value = 0.7501115169983111
$ws.Cells.Item(5, 35).Value = 0.7700900208196997
$ws.Cells.Item(5, 36).Value = 0.7226185418494169
$ws.Cells.Item(5, 37).Value = 0.7440450745358979
$ws.Cells.Item(5, 38).Value = 0
$ws.Cells.Item(5, 39).Value = 0
$ws.Cells.Item(5, 40).ClearContents()
$ws.Cells.Item(5, 42).ClearContents()

# Row 6
$ws.Cells.Item(6, 4).Value = -0.0106
$ws.Cells.Item(6, 5).Value = 0.0125
$ws.Cells.Item(6, 6).Value = 0.0775
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 6651.5
$ws.Cells.Item(6, 12).Value = 0.4480935057935866
$ws.Cells.Item(6, 13).Value = 5312.099999999999
$ws.Cells.Item(6, 14).Value = 0.04732366209830878
$ws.Cells.Item(6, 15).Value = 0.7986318875441629
$ws.Cells.Item(6, 16).Value = 5267.2
$ws.Cells.Item(6, 17).Value = 0.04692366352369346
$ws.Cells.Item(6, 18).Value = 0.7918815304818462
$ws.Cells.Item(6, 19).Value = 44.89999999999964
$ws.Cells.Item(6, 20).Value = 0.008452401121966763
$ws.Cells.Item(6, 21).Value = 18853.2
$ws.Cells.Item(6, 22).Value = 0.1679566397981656
$ws.Cells.Item(6, 23).Value = 0.136310167879861
$ws.Cells.Item(6, 24).Value = 0.05513041439641753
$ws.Cells.Item(6, 25).Value = 0.08117975348344345
$ws.Cells.Item(6, 26).Value = 0.07059639303985385
$ws.Cells.Item(6, 27).Value = 0
$ws.Cells.Item(6, 28).Value = 0.03018330337425221
$ws.Cells.Item(6, 29).Value = -0.03018330337425221
$ws.Cells.Item(6, 30).Value = 164499.4
$ws.Cells.Item(6, 31).Value = 0
$ws.Cells.Item(6, 32).Value = 164499.4
$ws.Cells.Item(6, 33).Value = 145646.2
$ws.Cells.Item(6, 34).Value = 0.59439753886001
$ws.Cells.Item(6, 35).Value = 0.7679050912526655
$ws.Cells.Item(6, 36).Value = 0.5647464914233068
$ws.Cells.Item(6, 37).Value = 0.7455073882144824
$ws.Cells.Item(6, 38).Value = 0
$ws.Cells.Item(6, 39).Value = 0
$ws.Cells.Item(6, 40).ClearContents()
$ws.Cells.Item(6, 42).ClearContents()
